$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Platform Coverage" sheet

# The coverage columns used throughout this sheet (every other column,
# corresponding to years 2026, 2028, 2030, ... 2040).
$covCols = @("P", "R", "T", "V", "X", "Z", "AB", "AD")

# --- Insert a new MDA age-band row at position 3 -------------------------
# This shifts the old rows 3-7 down to rows 4-8 (and their formatting,
# e.g. the numeric-text style on F/G of the EPI row, moves with them).
$ws.Rows.Item(3).Insert()

# --- Row 2 (All / Treatment / Campaign / MDA, ages 5-15) ------------------
# Drop the repeated 0.6 coverage values that used to run all the way to
# column AD; only the first four (H, J, L, N) remain.
$ws.Range("P2:AD2").ClearContents()

# --- Row 3 (new row: All / Treatment / Campaign / MDA, ages 2-15) ---------
$ws.Range("A3").Value = "All"
$ws.Range("B3").Value = "Treatment"
$ws.Range("C3").Value = "Campaign"
$ws.Range("D3").Value = "MDA"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 15
foreach ($col in $covCols) {
    $ws.Range($col + "3").Value = 0.8
}

# --- Row 4 (formerly row 3: All / Treatment / Campaign / MDA, ages 15-50) -
foreach ($col in $covCols) {
    $ws.Range($col + "4").Value = 0.5
}

# --- Row 5 (formerly row 4: All / Treatment / Campaign / MDA, ages 50-65) -
foreach ($col in $covCols) {
    $ws.Range($col + "5").Value = 0.5
}

# Rows 6-8 (formerly rows 5-7: the Vaccine/EPI, Vaccine/School and
# Vaccine/Out-of-school campaign rows) are left exactly as they were
# shifted to by the row insert above - no further changes needed there.

# --- Restore the sheet view (zoom level / selected cell) ------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("AD2").Select()
